$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort the TRANSACTIONS block (E3:H12) by the date column (H) ascending,
# keeping the CATEGORY block (A:C) untouched.
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("H3:H12"))
$sortObj.SetRange($ws.Range("E3:H12"))
$sortObj.Header = 2
$sortObj.Apply()

# Remove the TRANSACDATE table (columns J:K) entirely.
$ws.Range("J1:K12").EntireColumn.Delete()

# Restore a plausible selection position.
$ws.Range("G20").Select()
